$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text looks like a percentage need to be forced to Text format,
# otherwise Excel auto-converts "35%" into the number 0.35.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"

$ws.Range("D2").Value = "35%"
$ws.Range("E3").Value = "53%"
$ws.Range("C4").Value = "22°"
$ws.Range("B6").Value = "22°"
$ws.Range("E7").Value = "86%"
$ws.Range("E8").Value = "67%"

# Restore default cell style now that the text value is safely stored,
# so no visible style/number-format change is left behind.
$ws.Range("D2").Style = "Normal"
$ws.Range("E3").Style = "Normal"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Style = "Normal"
